$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the author name (A2, merged A2:B2)
$ws.Range("A2").Value = "Benoît Schopfer"

# Row 5 (existing entry) - shorten text, adjust hours
$ws.Range("B5").Value = "Constitution d'un groupe et discussion pour trouver une idée de projet."
$ws.Range("C5").Value = 1

# Row 6 (existing entry) - shorten text, adjust hours
$ws.Range("B6").Value = "Relecture et compétion du document de description du projet"
$ws.Range("C6").Value = 1

# Row 7 (new entry)
$ws.Range("A7").Value = 43157
$ws.Range("B7").Value = "Retour et échanges  avec le professeur à propos de notre proposition de projet et organisation entre nous."
$ws.Range("C7").Value = 1.5

# Row 8 (new entry)
$ws.Range("A8").Value = 43158
$ws.Range("B8").Value = "Spécifications de notre projet, définition des souhaits de chacun, choix des outils à implémenter, mockup, schéma de Dépendances Fonctionnelles afin d'organiser l'ordre d'implémentation de chaque fonctionnalité et début de la rédaction du cahier des charges."
$ws.Range("C8").Value = 5

# Keep the active selection roughly where the author last left it
$ws.Range("B11").Select()
